$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BEFORE ---
#   Row1 (all style s="1"): A1=id  B1=name  C1=description  D1=langCode  E1=isActive
#   Row2: A2=10007(n)  B2="Pre-Reg"  C2=description(str)  D2="eng"  E2=isActive(bool)
#
# --- AFTER ---
#   Row1 (style s="1", A1 removed entirely): B1=lang_code  C1=id  D1=name  E1=descr  F1=is_active
#   Rows2-13: A=index 0..11 (style s="1", number)  B=lang_code  C=id(number)
#             D=name  E=descr  F=is_active(bool, always TRUE)

# Grab the bold/centered/bordered header format (currently on A1) into the
# clipboard so it can be stamped onto the new header row and the new index
# column further down.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:F1").PasteSpecial(-4122) | Out-Null

# Wipe the old data row's contents (format for B2:E2 was already default).
$ws.Range("A2:E2").ClearContents() | Out-Null

# A1 disappears completely in the target (no header in col A anymore).
$ws.Range("A1").Clear() | Out-Null

# New header text.
$ws.Range("B1").Value = "lang_code"
$ws.Range("C1").Value = "id"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "descr"
$ws.Range("F1").Value = "is_active"

$data = @(
    @(0,  "eng", 10001, "Pre-Registration",       "Web portal for pre-registrations", $true),
    @(1,  "eng", 10002, "Registration Client",    "Desktop application for Registrations", $true),
    @(2,  "eng", 10003, "Registration Processor", "Application for post-registration process", $true),
    @(3,  "eng", 10004, "ID Authentication",      "Application for third party service provider authentication", $true),
    @(4,  "eng", 10005, "ID Control",              "Web portal for configuring applications", $true),
    @(5,  "eng", 10006, "Resident Portal",         "Web portal for Post ID generation services", $true),
    @(6,  "fra", 10001, "Pré-inscription",         "Portail Web pour les pré-inscriptions", $true),
    @(7,  "fra", 10002, "Client dinscription",     "Application de bureau pour les inscriptions", $true),
    @(8,  "fra", 10003, "Processeur dinscription", "Demande de post-inscription", $true),
    @(9,  "fra", 10004, "Authentification ID",     "Application pour lauthentification du fournisseur de services tiers", $true),
    @(10, "fra", 10005, "Contrôle didentité",      "Portail Web pour la configuration dapplications", $true),
    @(11, "fra", 10006, "Portail Résident",         "Portail Web pour les services de génération de post-ID", $true)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
    $row++
}

# Stamp the same header style onto the new index column A2:A13 (F1 still
# carries the header format we painted on earlier).
$ws.Range("F1").Copy() | Out-Null
$ws.Range("A2:A13").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
